# feat: add 2022-Q4 data
#
# "总计" (summary) sheet: insert the new 2022-Q4 totals as the new row 2,
# pushing the previous 2022-Q3 totals row down to row 3.
#
# "2022-Q3" sheet: duplicated so the original quarterly detail data is kept
# intact on a new sheet (after the copy operation, still named "2022-Q3"),
# while the original sheet object is renamed "2022-Q4" and its contents are
# replaced with the new quarter's fund holdings table.

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1) Update the "总计" summary sheet.
# ---------------------------------------------------------------------

# Push the existing 2022-Q3 summary row (row 2) down to row 3, keeping its
# formatting (column A uses the bold/bordered style already used by A2).
$summary.Range("A3").Value = 1
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.01

# Write the new 2022-Q4 totals into row 2.
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 6
$summary.Range("D2").Value = 0.04

# ---------------------------------------------------------------------
# 2) Duplicate the "2022-Q3" sheet so its data survives on its own sheet,
#    then turn the original sheet into the new "2022-Q4" sheet.
# ---------------------------------------------------------------------

$q3.Copy($null, $q3)
$q3copy = $wb.Worksheets.Item(3)
$q3.Name = "2022-Q4"
$q3copy.Name = "2022-Q3"

$q4 = $q3

# Clear all prior contents/formatting of the (renamed) Q4 sheet so it can
# be rebuilt from scratch with the new fund table.
$q4.Cells.Clear()

# Header row.
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Fund holdings rows.
$rows = @(
    @("010857", "宝盈祥乐一年持有期混合型证券投资基金A", "1.00", "33.11", "1.25", "0.0125", 10),
    @("014246", "大摩现代服务业混合A",                     "0.17", "86.98", "5.89", "0.0100", 8),
    @("008324", "宝盈祥利稳健配置混合A",                   "0.53", "31.08", "1.25", "0.0066", 7),
    @("008325", "宝盈祥利稳健配置混合C",                   "0.29", "31.08", "1.25", "0.0036", 7),
    @("014247", "大摩现代服务业混合C",                     "0.06", "86.98", "5.89", "0.0035", 8),
    @("010858", "宝盈祥乐一年持有期混合型证券投资基金C",   "0.06", "33.11", "1.25", "0.0008", 10)
)

# Column B (fund codes, which have significant leading zeros) and columns
# D:G (text-looking numbers, e.g. "1.00") are stored as text, just like the
# source data, so force text formatting before writing them.
$q4.Range("B2:B7").NumberFormat = "@"
$q4.Range("D2:G7").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $q4.Cells.Item($r, 1).Value = $i
    $q4.Cells.Item($r, 2).Value = $row[0]
    $q4.Cells.Item($r, 3).Value = $row[1]
    $q4.Cells.Item($r, 4).Value = $row[2]
    $q4.Cells.Item($r, 5).Value = $row[3]
    $q4.Cells.Item($r, 6).Value = $row[4]
    $q4.Cells.Item($r, 7).Value = $row[5]
    $q4.Cells.Item($r, 8).Value = $row[6]
}

# Drop the temporary text number-format now that the values are written,
# matching the unstyled cells in the target sheet.
$q4.Range("B2:B7").Style = "Normal"
$q4.Range("D2:G7").Style = "Normal"

# Apply the bold/bordered header-and-index-column style (the same style
# already used on the "总计" sheet) to the new header row and index column.
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$summary.Range("A2").Copy()
$q4.Range("A2:A7").PasteSpecial(-4122)

# The sheet selection itself is unchanged by this edit (the "总计" sheet
# was, and remains, the active tab in the workbook view).
$summary.Activate()
